$d = $word.ActiveDocument

# Locate the first "[contributions here]" placeholder -- the one that
# follows "Cheah Meng Yew" -- and turn it into the first of three
# bullet items describing his contribution.
$rng = $d.Content
$rng.Find.Execute("[contributions here]", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "LETTER/NUMBER SEPERATION coding contribution "

# Add a second list item (same list formatting) right after it.
$lq = [char]0x201C
$rq = [char]0x201D
$line2 = "DOCUMENTATION on " + $lq + "Overall Methodology" + $rq + " and " + $lq + "Detailed Description" + $rq + " "
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null
$rng.InsertAfter($line2)

# Add a third list item (same list formatting) right after that one.
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null
$rng.InsertAfter("SEARCHING AND TESTING on additional Car Plates ")
